$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:19:40.572088"
$ws1.Range("F3").Value = "2021-10-05 14:19:40.572096"
$ws1.Range("F4").Value = "2021-10-05 14:19:40.572099"
$ws1.Range("F5").Value = "2021-10-05 14:19:40.572101"
$ws1.Range("F6").Value = "2021-10-05 14:19:40.572104"
$ws1.Range("F7").Value = "2021-10-05 14:19:40.572107"
$ws1.Range("F8").Value = "2021-10-05 14:19:40.572110"
$ws1.Range("F9").Value = "2021-10-05 14:19:40.572112"
$ws1.Range("F10").Value = "2021-10-05 14:19:40.572115"
$ws1.Range("F11").Value = "2021-10-05 14:19:40.572117"
$ws1.Range("F12").Value = "2021-10-05 14:19:40.572120"
$ws1.Range("F13").Value = "2021-10-05 14:19:40.572123"
$ws1.Range("F14").Value = "2021-10-05 14:19:40.572125"
$ws1.Range("F15").Value = "2021-10-05 14:19:40.572128"
$ws1.Range("F16").Value = "2021-10-05 14:19:40.572130"
$ws1.Range("F17").Value = "2021-10-05 14:19:40.572132"
$ws1.Range("F18").Value = "2021-10-05 14:19:40.572135"
$ws1.Range("F19").Value = "2021-10-05 14:19:40.572138"
$ws1.Range("F20").Value = "2021-10-05 14:19:40.572140"
$ws1.Range("F21").Value = "2021-10-05 14:19:40.572143"
$ws1.Range("F22").Value = "2021-10-05 14:19:40.572145"
$ws1.Range("F23").Value = "2021-10-05 14:19:40.572148"

# Add metadata worksheet after the "data" sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row (bold/bordered style copied from the data sheet's header cells)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

# Data row
$ws2.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B2").Value = "Confirmed Fanconi anaemia or Bloom syndrome"
$ws2.Range("C2").Value = 508
$ws2.Range("D2").Value = "'1.11"
$ws2.Range("D2").ClearFormats()
$ws2.Range("E2").Value = "2020-12-02T15:34:58.834861Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:40.568219"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/508/?format=json"

# Keep "data" as the active sheet/tab, matching the unchanged bookViews in the diff
$ws1.Activate()
